$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update calibration values for L3 leg (row 6)
$ws.Range("C6").Value = 1350
$ws.Range("D6").Value = 1550

# Update calibration value for R2 leg (row 8)
$ws.Range("D8").Value = 1250

# Update the active cell selection shown in the sheet view
$ws.Range("E17").Select()
